$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 2.317355952907718 / 1000000
$ws.Range("C2").Value = 7.097389502863649 / 100000
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 671.5447565614272

$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 10.29869402782916
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 24.31896926437656
